# Add the missing full stop ('.') at the end of the two list items that
# are currently missing their sentence-terminating period:
#   1) "මහා මාර්ග ප්‍රතිසංස්කරණය"
#   2) "අත්තුපනායික ධර්ම පරියාය"
# The new run uses the "Nirmala UI" font for ascii/hAnsi/cs, matching the
# formatting already used by the surrounding Sinhala-text runs.

$d = $word.ActiveDocument

function Add-TrailingPeriod([string]$searchText) {
    # Search the whole document content for the anchor text. Running
    # Find.Execute directly on a Range collapses/repositions that Range to
    # the found match (standard Word Range.Find semantics).
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $false, $false, $false, $false,
                                $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return
    }

    # Grow the matched range out to the full enclosing paragraph, then
    # collapse to just before the paragraph mark (the last character of a
    # paragraph Range is always the pilcrow), so the new run lands at the
    # very end of the paragraph's text.
    $rng.Expand(4)
    $rng.SetRange($rng.End - 1, $rng.End - 1)
    $rng.InsertAfter(".")

    # Apply the "Nirmala UI" font (ascii/hAnsi via Font.Name, complex
    # script via Font.NameBi) to just the newly inserted period, matching
    # the formatting used by the neighboring Sinhala-text runs.
    $find = $rng.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Replacement.Font.Name = "Nirmala UI"
    $find.Replacement.Font.NameBi = "Nirmala UI"
    $find.Execute(".", $false, $false, $false, $false, $false, $true, 1,
                  $false, ".", 2)
}

Add-TrailingPeriod("රතිසංස්කරණය")
Add-TrailingPeriod("පරියාය")
